{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright paragraph block\n// (and the blank paragraph immediately before it) that used to follow\n// the \"LOQ4009: Instrumenta\u00e7\u00e3o na Industria Qu\u00edmica (Requisito fraco)\"\n// requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"LOQ4009: Instrumenta\u00e7\u00e3o na Industria Qu\u00edmica (Requisito fraco)\";\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1 && markerIndex + 3 < items.length) {\n  // The three paragraphs immediately following the marker are, in order:\n  // an empty paragraph, \"Ver no Jupiter Salvar em pdf Salvar em docx\", and\n  // the \"\u00a9 2020 ...\" copyright line. Verify before deleting so the script\n  // fails loudly instead of silently removing the wrong content.\n  const blank = items[markerIndex + 1].text.trim();\n  const first = items[markerIndex + 2].text.trim();\n  const second = items[markerIndex + 3].text.trim();\n\n  if (blank === \"\" && first === targets[0] && second === targets[1]) {\n    items[markerIndex + 3].delete();\n    items[markerIndex + 2].delete();\n    items[markerIndex + 1].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright paragraph block\n# (and the blank paragraph immediately before it) that used to follow\n# the \"LOQ4009: Instrumenta\u00e7\u00e3o na Industria Qu\u00edmica (Requisito fraco)\"\n# requirement line.\n$d = $word.ActiveDocument\n\n$marker = \"LOQ4009: Instrumenta\u00e7\u00e3o na Industria Qu\u00edmica (Requisito fraco)\"\n$firstTarget = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$secondTarget = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq $marker) {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -ne -1 -and ($markerIndex + 3) -le $d.Paragraphs.Count) {\n    $blank = $d.Paragraphs.Item($markerIndex + 1).Range.Text.Trim()\n    $first = $d.Paragraphs.Item($markerIndex + 2).Range.Text.Trim()\n    $second = $d.Paragraphs.Item($markerIndex + 3).Range.Text.Trim()\n\n    if ($blank -eq \"\" -and $first -eq $firstTarget -and $second -eq $secondTarget) {\n        $startPara = $d.Paragraphs.Item($markerIndex + 1)\n        $endPara = $d.Paragraphs.Item($markerIndex + 3)\n        $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n        $range.Delete()\n    }\n}\n"}
